# Apply scraped-schedule update for Línea 141 workbook.
# New scrape timestamp: 07:17:57 (previous: 06:52:31)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("LP1912")
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws3 = $wb.Worksheets.Item("6203-6173")

# --- Update header info rows (A2 = last update time, A3 = total row count) ---
$ws1.Cells.Item(2, 1).Value = "Última actualización: 07:17:57"
$ws1.Cells.Item(3, 1).Value = "Total filas: 79"

$ws2.Cells.Item(2, 1).Value = "Última actualización: 07:17:57"
$ws2.Cells.Item(3, 1).Value = "Total filas: 19"

$ws3.Cells.Item(2, 1).Value = "Última actualización: 07:17:57"
$ws3.Cells.Item(3, 1).Value = "Total filas: 18"

# --- Sheet1 (LP1912) data rows 51-84 ---
$ws1.Cells.Item(51, 1).Value = "07:17:57"
$ws1.Cells.Item(51, 2).Value = "07:20"
$ws1.Cells.Item(51, 3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(51, 4).Value = 3
$ws1.Cells.Item(51, 5).Value = "LP1912"
$ws1.Cells.Item(52, 1).Value = "06:52:31"
$ws1.Cells.Item(52, 2).Value = "07:21"
$ws1.Cells.Item(52, 3).Value = "10_OLMOS"
$ws1.Cells.Item(52, 4).Value = 29
$ws1.Cells.Item(52, 5).Value = "LP1912"
$ws1.Cells.Item(53, 1).Value = "07:17:57"
$ws1.Cells.Item(53, 2).Value = "07:22"
$ws1.Cells.Item(53, 3).Value = "10_OLMOS"
$ws1.Cells.Item(53, 4).Value = 5
$ws1.Cells.Item(53, 5).Value = "LP1912"
$ws1.Cells.Item(54, 1).Value = "06:52:31"
$ws1.Cells.Item(54, 2).Value = "07:23"
$ws1.Cells.Item(54, 3).Value = "10_OLMOS"
$ws1.Cells.Item(54, 4).Value = 31
$ws1.Cells.Item(54, 5).Value = "LP1912"
$ws1.Cells.Item(55, 1).Value = "07:17:57"
$ws1.Cells.Item(55, 2).Value = "07:31"
$ws1.Cells.Item(55, 3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(55, 4).Value = 14
$ws1.Cells.Item(55, 5).Value = "LP1912"
$ws1.Cells.Item(56, 1).Value = "07:17:57"
$ws1.Cells.Item(56, 2).Value = "07:31"
$ws1.Cells.Item(56, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(56, 4).Value = 14
$ws1.Cells.Item(56, 5).Value = "LP1912"
$ws1.Cells.Item(57, 1).Value = "07:17:57"
$ws1.Cells.Item(57, 2).Value = "07:31"
$ws1.Cells.Item(57, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(57, 4).Value = 14
$ws1.Cells.Item(57, 5).Value = "LP1912"
$ws1.Cells.Item(58, 1).Value = "06:52:31"
$ws1.Cells.Item(58, 2).Value = "07:32"
$ws1.Cells.Item(58, 3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(58, 4).Value = 40
$ws1.Cells.Item(58, 5).Value = "LP1912"
$ws1.Cells.Item(59, 1).Value = "07:17:57"
$ws1.Cells.Item(59, 2).Value = "07:35"
$ws1.Cells.Item(59, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(59, 4).Value = 18
$ws1.Cells.Item(59, 5).Value = "LP1912"
$ws1.Cells.Item(60, 1).Value = "07:17:57"
$ws1.Cells.Item(60, 2).Value = "07:36"
$ws1.Cells.Item(60, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(60, 4).Value = 19
$ws1.Cells.Item(60, 5).Value = "LP1912"
$ws1.Cells.Item(61, 1).Value = "07:17:57"
$ws1.Cells.Item(61, 2).Value = "07:38"
$ws1.Cells.Item(61, 3).Value = "10_OLMOS"
$ws1.Cells.Item(61, 4).Value = 21
$ws1.Cells.Item(61, 5).Value = "LP1912"
$ws1.Cells.Item(62, 1).Value = "06:52:31"
$ws1.Cells.Item(62, 2).Value = "07:39"
$ws1.Cells.Item(62, 3).Value = "10_OLMOS"
$ws1.Cells.Item(62, 4).Value = 47
$ws1.Cells.Item(62, 5).Value = "LP1912"
$ws1.Cells.Item(63, 1).Value = "07:17:57"
$ws1.Cells.Item(63, 2).Value = "07:46"
$ws1.Cells.Item(63, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(63, 4).Value = 29
$ws1.Cells.Item(63, 5).Value = "LP1912"
$ws1.Cells.Item(64, 1).Value = "07:17:57"
$ws1.Cells.Item(64, 2).Value = "07:47"
$ws1.Cells.Item(64, 3).Value = "14_ABASTO"
$ws1.Cells.Item(64, 4).Value = 30
$ws1.Cells.Item(64, 5).Value = "LP1912"
$ws1.Cells.Item(65, 1).Value = "07:17:57"
$ws1.Cells.Item(65, 2).Value = "07:51"
$ws1.Cells.Item(65, 3).Value = "215D_EL PATO"
$ws1.Cells.Item(65, 4).Value = 34
$ws1.Cells.Item(65, 5).Value = "LP1912"
$ws1.Cells.Item(66, 1).Value = "07:17:57"
$ws1.Cells.Item(66, 2).Value = "07:59"
$ws1.Cells.Item(66, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(66, 4).Value = 42
$ws1.Cells.Item(66, 5).Value = "LP1912"
$ws1.Cells.Item(67, 1).Value = "07:17:57"
$ws1.Cells.Item(67, 2).Value = "08:03"
$ws1.Cells.Item(67, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(67, 4).Value = 46
$ws1.Cells.Item(67, 5).Value = "LP1912"
$ws1.Cells.Item(68, 1).Value = "07:17:57"
$ws1.Cells.Item(68, 2).Value = "08:11"
$ws1.Cells.Item(68, 3).Value = "15_ABASTO"
$ws1.Cells.Item(68, 4).Value = 54
$ws1.Cells.Item(68, 5).Value = "LP1912"
$ws1.Cells.Item(69, 1).Value = "06:52:31"
$ws1.Cells.Item(69, 2).Value = "08:12"
$ws1.Cells.Item(69, 3).Value = "15_ABASTO"
$ws1.Cells.Item(69, 4).Value = 80
$ws1.Cells.Item(69, 5).Value = "LP1912"
$ws1.Cells.Item(70, 1).Value = "07:17:57"
$ws1.Cells.Item(70, 2).Value = "08:20"
$ws1.Cells.Item(70, 3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(70, 4).Value = 63
$ws1.Cells.Item(70, 5).Value = "LP1912"
$ws1.Cells.Item(71, 1).Value = "06:52:31"
$ws1.Cells.Item(71, 2).Value = "08:21"
$ws1.Cells.Item(71, 3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(71, 4).Value = 89
$ws1.Cells.Item(71, 5).Value = "LP1912"
$ws1.Cells.Item(72, 1).Value = "07:17:57"
$ws1.Cells.Item(72, 2).Value = "08:22"
$ws1.Cells.Item(72, 3).Value = "215B_EL PATO"
$ws1.Cells.Item(72, 4).Value = 65
$ws1.Cells.Item(72, 5).Value = "LP1912"
$ws1.Cells.Item(73, 1).Value = "07:17:57"
$ws1.Cells.Item(73, 2).Value = "08:22"
$ws1.Cells.Item(73, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(73, 4).Value = 65
$ws1.Cells.Item(73, 5).Value = "LP1912"
$ws1.Cells.Item(74, 1).Value = "06:52:31"
$ws1.Cells.Item(74, 2).Value = "08:23"
$ws1.Cells.Item(74, 3).Value = "215B_EL PATO"
$ws1.Cells.Item(74, 4).Value = 91
$ws1.Cells.Item(74, 5).Value = "LP1912"
$ws1.Cells.Item(75, 1).Value = "07:17:57"
$ws1.Cells.Item(75, 2).Value = "08:26"
$ws1.Cells.Item(75, 3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(75, 4).Value = 69
$ws1.Cells.Item(75, 5).Value = "LP1912"
$ws1.Cells.Item(76, 1).Value = "06:52:31"
$ws1.Cells.Item(76, 2).Value = "08:27"
$ws1.Cells.Item(76, 3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(76, 4).Value = 95
$ws1.Cells.Item(76, 5).Value = "LP1912"
$ws1.Cells.Item(77, 1).Value = "06:52:31"
$ws1.Cells.Item(77, 2).Value = "08:35"
$ws1.Cells.Item(77, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(77, 4).Value = 103
$ws1.Cells.Item(77, 5).Value = "LP1912"
$ws1.Cells.Item(78, 1).Value = "07:17:57"
$ws1.Cells.Item(78, 2).Value = "08:41"
$ws1.Cells.Item(78, 3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(78, 4).Value = 84
$ws1.Cells.Item(78, 5).Value = "LP1912"
$ws1.Cells.Item(79, 1).Value = "06:52:31"
$ws1.Cells.Item(79, 2).Value = "08:42"
$ws1.Cells.Item(79, 3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(79, 4).Value = 110
$ws1.Cells.Item(79, 5).Value = "LP1912"
$ws1.Cells.Item(80, 1).Value = "07:17:57"
$ws1.Cells.Item(80, 2).Value = "08:43"
$ws1.Cells.Item(80, 3).Value = "14_ABASTO"
$ws1.Cells.Item(80, 4).Value = 86
$ws1.Cells.Item(80, 5).Value = "LP1912"
$ws1.Cells.Item(81, 1).Value = "07:17:57"
$ws1.Cells.Item(81, 2).Value = "08:53"
$ws1.Cells.Item(81, 3).Value = "17_ROMERO"
$ws1.Cells.Item(81, 4).Value = 96
$ws1.Cells.Item(81, 5).Value = "LP1912"
$ws1.Cells.Item(82, 1).Value = "07:17:57"
$ws1.Cells.Item(82, 2).Value = "09:01"
$ws1.Cells.Item(82, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(82, 4).Value = 104
$ws1.Cells.Item(82, 5).Value = "LP1912"
$ws1.Cells.Item(83, 1).Value = "07:17:57"
$ws1.Cells.Item(83, 2).Value = "09:10"
$ws1.Cells.Item(83, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(83, 4).Value = 113
$ws1.Cells.Item(83, 5).Value = "LP1912"
$ws1.Cells.Item(84, 1).Value = "07:17:57"
$ws1.Cells.Item(84, 2).Value = "09:16"
$ws1.Cells.Item(84, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(84, 4).Value = 119
$ws1.Cells.Item(84, 5).Value = "LP1912"

# --- Sheet2 (LP1912-215) data rows 21-24 ---
$ws2.Cells.Item(21, 1).Value = "07:17:57"
$ws2.Cells.Item(21, 2).Value = "07:51"
$ws2.Cells.Item(21, 3).Value = "215D_EL PATO"
$ws2.Cells.Item(21, 4).Value = 34
$ws2.Cells.Item(21, 5).Value = "LP1912"
$ws2.Cells.Item(22, 1).Value = "07:17:57"
$ws2.Cells.Item(22, 2).Value = "08:22"
$ws2.Cells.Item(22, 3).Value = "215B_EL PATO"
$ws2.Cells.Item(22, 4).Value = 65
$ws2.Cells.Item(22, 5).Value = "LP1912"
$ws2.Cells.Item(23, 1).Value = "06:52:31"
$ws2.Cells.Item(23, 2).Value = "08:23"
$ws2.Cells.Item(23, 3).Value = "215B_EL PATO"
$ws2.Cells.Item(23, 4).Value = 91
$ws2.Cells.Item(23, 5).Value = "LP1912"
$ws2.Cells.Item(24, 1).Value = "07:17:57"
$ws2.Cells.Item(24, 2).Value = "09:01"
$ws2.Cells.Item(24, 3).Value = "215A_EL PATO"
$ws2.Cells.Item(24, 4).Value = 104
$ws2.Cells.Item(24, 5).Value = "LP1912"

# --- Sheet3 (6203-6173) data rows 18-23 ---
$ws3.Cells.Item(18, 1).Value = "07:17:57"
$ws3.Cells.Item(18, 2).Value = "07:48"
$ws3.Cells.Item(18, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(18, 4).Value = 31
$ws3.Cells.Item(18, 5).Value = "L6173"
$ws3.Cells.Item(19, 1).Value = "06:52:31"
$ws3.Cells.Item(19, 2).Value = "08:07"
$ws3.Cells.Item(19, 3).Value = "215C_LA PLATA"
$ws3.Cells.Item(19, 4).Value = 75
$ws3.Cells.Item(19, 5).Value = "L6203"
$ws3.Cells.Item(20, 1).Value = "07:17:57"
$ws3.Cells.Item(20, 2).Value = "08:09"
$ws3.Cells.Item(20, 3).Value = "215C_LA PLATA"
$ws3.Cells.Item(20, 4).Value = 52
$ws3.Cells.Item(20, 5).Value = "L6203"
$ws3.Cells.Item(21, 1).Value = "06:52:31"
$ws3.Cells.Item(21, 2).Value = "08:30"
$ws3.Cells.Item(21, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(21, 4).Value = 98
$ws3.Cells.Item(21, 5).Value = "L6173"
$ws3.Cells.Item(22, 1).Value = "07:17:57"
$ws3.Cells.Item(22, 2).Value = "08:34"
$ws3.Cells.Item(22, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(22, 4).Value = 77
$ws3.Cells.Item(22, 5).Value = "L6173"
$ws3.Cells.Item(23, 1).Value = "07:17:57"
$ws3.Cells.Item(23, 2).Value = "09:08"
$ws3.Cells.Item(23, 3).Value = "215D_LA PLATA"
$ws3.Cells.Item(23, 4).Value = 111
$ws3.Cells.Item(23, 5).Value = "L6203"
